$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add Sheet2 right after Sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Type]::Missing, $ws1)
$ws2.Name = "Sheet2"

# --- Add Sheet3 right after Sheet2 ---
$ws3 = $wb.Worksheets.Add([System.Type]::Missing, $ws2)
$ws3.Name = "Sheet3"

# --- Duplicate Sheet1's data (A1:C8 - headers x/y/value + 7 data rows) onto Sheet2 and Sheet3 ---
for ($r = 1; $r -le 8; $r++) {
    for ($c = 1; $c -le 3; $c++) {
        $val = $ws1.Cells.Item($r, $c).Value()
        $ws2.Cells.Item($r, $c).Value = $val
        $ws3.Cells.Item($r, $c).Value = $val
    }
}

# --- Restore Sheet1's selection to the full data range ---
$ws1.Activate() | Out-Null
$ws1.Range("A1:C8").Select() | Out-Null

# --- Sheet2: full data range selected as well ---
$ws2.Activate() | Out-Null
$ws2.Range("A1:C8").Select() | Out-Null

# --- Sheet3 ends up the active tab, with a single cell (J9) selected ---
$ws3.Activate() | Out-Null
$ws3.Range("J9").Select() | Out-Null
